$wb = $excel.ActiveWorkbook

# Overview sheet: Latest HO Xliff Generate Date for d42de85d... row
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G4").Value = "2016-08-23 20:47:48"

# zh-cn sheet: Correspond Handoff Datetime / Correspond Handback DateTime for d42de85d... row
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H4").Value = "2016-08-23 20:47:44"
$wsZhCn.Range("K4").Value = "2016-08-23 20:48:05"

# de-de sheet: Correspond Handback DateTime for d42de85d... row
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("K4").Value = "2016-08-23 20:48:17"
